# Templ.docx edit: replace the legacy "Table1.<Field>" merge placeholders
# in the first table's second row with the new "[Questions]<Field>" style
# placeholders (and expand the "Theme" cell into the six Questions-scoped
# fields that replace it).

$d = $word.ActiveDocument
$t1 = $d.Tables.Item(1)

# Column 1, row 2: <Table1.NumTheme> -> <[Questions]Sequence>
$cell1 = $t1.Cell(2, 1)
$cell1.Range.Find.Execute("<Table1.NumTheme>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<[Questions]Sequence>", 2)

# Column 2, row 2: <Table1.Theme> -> six chained <[Questions]...> placeholders
$cell2 = $t1.Cell(2, 2)
$cell2.Range.Find.Execute("<Table1.Theme>", $false, $false, $false, $false, $false, `
    $true, 1, $false, `
    "<[Questions]Subject> <[Questions]Direction> <[Questions]Tarif> <[Questions]TarifView> <[Questions]TechConnection> <[Questions]Organization>", `
    2)

# Column 3, row 2: <Table1.AuthPerson> -> <[Questions]AuthPerson>
$cell3 = $t1.Cell(2, 3)
$cell3.Range.Find.Execute("<Table1.AuthPerson>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<[Questions]AuthPerson>", 2)
